# DATA_goal/Junction_Flooding_447.xlsx edit
# - refresh the J1..J33 flooding sample rows with a new 4-row slice of data (1000-row dataset)
# - drop the now-unused trailing 6th sample row
# - nudge a batch of data columns one character wider ("custom accuracy" formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- replace the sampled sensor readings in rows 2-5 ----
# row 2 (column A = timestamp, B..AH = J1..J33 readings)
$ws.Cells.Item(2, 1).Value = 45128.50694444445
$ws.Cells.Item(2, 2).Value = 13.798
$ws.Cells.Item(2, 3).Value = 9.137
$ws.Cells.Item(2, 4).Value = 3.527
$ws.Cells.Item(2, 5).Value = 29.879
$ws.Cells.Item(2, 6).Value = 22.444
$ws.Cells.Item(2, 7).Value = 10.657
$ws.Cells.Item(2, 8).Value = 31.967
$ws.Cells.Item(2, 9).Value = 16.87
$ws.Cells.Item(2, 10).Value = 6.742
$ws.Cells.Item(2, 11).Value = 10.011
$ws.Cells.Item(2, 12).Value = 11.733
$ws.Cells.Item(2, 13).Value = 12.516
$ws.Cells.Item(2, 14).Value = 3.497
$ws.Cells.Item(2, 15).Value = 10.903
$ws.Cells.Item(2, 16).Value = 14.966
$ws.Cells.Item(2, 17).Value = 9.704000000000001
$ws.Cells.Item(2, 18).Value = 3.059
$ws.Cells.Item(2, 19).Value = 1.672
$ws.Cells.Item(2, 20).Value = 158.575
$ws.Cells.Item(2, 21).Value = 30.182
$ws.Cells.Item(2, 22).Value = 10.064
$ws.Cells.Item(2, 23).Value = 19.331
$ws.Cells.Item(2, 24).Value = 9.888999999999999
$ws.Cells.Item(2, 25).Value = 2.87
$ws.Cells.Item(2, 26).Value = 17.102
$ws.Cells.Item(2, 27).Value = 8.888999999999999
$ws.Cells.Item(2, 28).Value = 8.15
$ws.Cells.Item(2, 29).Value = 9.673999999999999
$ws.Cells.Item(2, 30).Value = 12.058
$ws.Cells.Item(2, 31).Value = 3.072
$ws.Cells.Item(2, 32).Value = 28.934
$ws.Cells.Item(2, 33).Value = 5.407
$ws.Cells.Item(2, 34).Value = 12.581

# row 3 (column A = timestamp, B..AH = J1..J33 readings)
$ws.Cells.Item(3, 1).Value = 45128.51388888889
$ws.Cells.Item(3, 2).Value = 22.468
$ws.Cells.Item(3, 3).Value = 16.378
$ws.Cells.Item(3, 4).Value = 1.963
$ws.Cells.Item(3, 5).Value = 49.095
$ws.Cells.Item(3, 6).Value = 39.376
$ws.Cells.Item(3, 7).Value = 17.56
$ws.Cells.Item(3, 8).Value = 65.952
$ws.Cells.Item(3, 9).Value = 27.341
$ws.Cells.Item(3, 10).Value = 11.927
$ws.Cells.Item(3, 11).Value = 17.56
$ws.Cells.Item(3, 12).Value = 19.612
$ws.Cells.Item(3, 13).Value = 20.857
$ws.Cells.Item(3, 14).Value = 5.676
$ws.Cells.Item(3, 15).Value = 17.67
$ws.Cells.Item(3, 16).Value = 24.975
$ws.Cells.Item(3, 17).Value = 15.139
$ws.Cells.Item(3, 18).Value = 1.502
$ws.Cells.Item(3, 19).Value = 1.201
$ws.Cells.Item(3, 20).Value = 261.645
$ws.Cells.Item(3, 21).Value = 49.388
$ws.Cells.Item(3, 22).Value = 16.31
$ws.Cells.Item(3, 23).Value = 32.863
$ws.Cells.Item(3, 24).Value = 17.177
$ws.Cells.Item(3, 25).Value = 2.938
$ws.Cells.Item(3, 26).Value = 32.822
$ws.Cells.Item(3, 27).Value = 14.407
$ws.Cells.Item(3, 28).Value = 12.891
$ws.Cells.Item(3, 29).Value = 15.184
$ws.Cells.Item(3, 30).Value = 20.441
$ws.Cells.Item(3, 31).Value = 1.198
$ws.Cells.Item(3, 32).Value = 60.346
$ws.Cells.Item(3, 33).Value = 9.081
$ws.Cells.Item(3, 34).Value = 20.392

# row 4 (column A = timestamp, B..AH = J1..J33 readings)
$ws.Cells.Item(4, 1).Value = 45128.52083333334
$ws.Cells.Item(4, 2).Value = 10.954
$ws.Cells.Item(4, 3).Value = 7.94
$ws.Cells.Item(4, 4).Value = 1.121
$ws.Cells.Item(4, 5).Value = 24.048
$ws.Cells.Item(4, 6).Value = 19.043
$ws.Cells.Item(4, 7).Value = 8.537000000000001
$ws.Cells.Item(4, 8).Value = 36.99
$ws.Cells.Item(4, 9).Value = 13.38
$ws.Cells.Item(4, 10).Value = 5.799
$ws.Cells.Item(4, 11).Value = 8.416
$ws.Cells.Item(4, 12).Value = 9.605
$ws.Cells.Item(4, 13).Value = 10.261
$ws.Cells.Item(4, 14).Value = 2.78
$ws.Cells.Item(4, 15).Value = 8.647
$ws.Cells.Item(4, 16).Value = 12.194
$ws.Cells.Item(4, 17).Value = 7.533
$ws.Cells.Item(4, 18).Value = 0.95
$ws.Cells.Item(4, 19).Value = 0.653
$ws.Cells.Item(4, 20).Value = 124.311
$ws.Cells.Item(4, 21).Value = 24.307
$ws.Cells.Item(4, 22).Value = 7.982
$ws.Cells.Item(4, 23).Value = 16.06
$ws.Cells.Item(4, 24).Value = 8.364000000000001
$ws.Cells.Item(4, 25).Value = 1.544
$ws.Cells.Item(4, 26).Value = 17.591
$ws.Cells.Item(4, 27).Value = 7.05
$ws.Cells.Item(4, 28).Value = 6.371
$ws.Cells.Item(4, 29).Value = 7.492
$ws.Cells.Item(4, 30).Value = 9.987
$ws.Cells.Item(4, 31).Value = 0.746
$ws.Cells.Item(4, 32).Value = 34.003
$ws.Cells.Item(4, 33).Value = 4.386
$ws.Cells.Item(4, 34).Value = 9.978999999999999

# row 5 (column A = timestamp, B..AH = J1..J33 readings)
$ws.Cells.Item(5, 1).Value = 45128.52777777778
$ws.Cells.Item(5, 2).Value = 14.81
$ws.Cells.Item(5, 3).Value = 10.93
$ws.Cells.Item(5, 4).Value = 1.04
$ws.Cells.Item(5, 5).Value = 32.41
$ws.Cells.Item(5, 6).Value = 26.16
$ws.Cells.Item(5, 7).Value = 11.59
$ws.Cells.Item(5, 8).Value = 45
$ws.Cells.Item(5, 9).Value = 18.03
$ws.Cells.Item(5, 10).Value = 7.93
$ws.Cells.Item(5, 11).Value = 11.66
$ws.Cells.Item(5, 12).Value = 12.98
$ws.Cells.Item(5, 13).Value = 13.82
$ws.Cells.Item(5, 14).Value = 3.74
$ws.Cells.Item(5, 15).Value = 11.65
$ws.Cells.Item(5, 16).Value = 16.51
$ws.Cells.Item(5, 17).Value = 9.94
$ws.Cells.Item(5, 18).Value = 0.77
$ws.Cells.Item(5, 19).Value = 0.67
$ws.Cells.Item(5, 20).Value = 170.07
$ws.Cells.Item(5, 21).Value = 32.57
$ws.Cells.Item(5, 22).Value = 10.76
$ws.Cells.Item(5, 23).Value = 21.75
$ws.Cells.Item(5, 24).Value = 11.41
$ws.Cells.Item(5, 25).Value = 1.84
$ws.Cells.Item(5, 26).Value = 21.87
$ws.Cells.Item(5, 27).Value = 9.5
$ws.Cells.Item(5, 28).Value = 8.48
$ws.Cells.Item(5, 29).Value = 9.970000000000001
$ws.Cells.Item(5, 30).Value = 13.57
$ws.Cells.Item(5, 31).Value = 0.54
$ws.Cells.Item(5, 32).Value = 40.9
$ws.Cells.Item(5, 33).Value = 6
$ws.Cells.Item(5, 34).Value = 13.45

# ---- the refreshed slice only has 4 data rows, so drop the old row 6 ----
$ws.Rows.Item(6).Delete()

# ---- widen a subset of columns by one character (ColumnWidth is offset from the
#      stored column <col width> by the standard ~0.83 char gridline/padding amount) ----
$colWidthPad = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(3).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(5).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(6).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(7).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(9).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(10).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(11).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(12).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(13).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(15).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(16).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(17).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(20).ColumnWidth = 9 - $colWidthPad
$ws.Columns.Item(21).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(22).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(23).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(24).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(26).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(27).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(28).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(29).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(30).ColumnWidth = 8 - $colWidthPad
$ws.Columns.Item(34).ColumnWidth = 8 - $colWidthPad
